$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin ranking table "Price" (D) and "Volume(1h)" (E) refresh.
# D-column values are textual price strings (e.g. "29.251.87", "1.000") that
# must stay literal text, not be re-interpreted as numbers, so each D cell
# that changes is switched to Text format before its value is written.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.251.87'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.844.93'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '242.86'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = '0.6624'
$ws.Range("E6").Value = '  -1.18%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '44.91'
$ws.Range("E8").Value = '  +6.69%  '
$ws.Range("D9").Value = '0.07445'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = '0.2958'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = '23.32'
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").Value = '1.842.66'
$ws.Range("E13").Value = '  +14.15%  '
$ws.Range("D14").Value = '5.025'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = '0.6731'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").Value = '83.44'
$ws.Range("E16").Value = '  -3.42%  '
$ws.Range("D17").Value = '6.183'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = '0.000008760'
$ws.Range("E18").Value = '  +5.89%  '
$ws.Range("D19").Value = '29.249.15'
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("D20").Value = '2.094.63'
$ws.Range("E20").Value = '  +3.73%  '
$ws.Range("D21").Value = '12.55'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = '227.04'
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '7.163'
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '158.98'
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '8.632'
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '0.1406'
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '1.510'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D32").Value = '4.064'
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("D34").Value = '0.05333'
$ws.Range("E34").Value = '  -0.66%  '
$ws.Range("D35").Value = '1.868'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("D36").Value = '0.7489'
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").Value = '1.313.10'
$ws.Range("E39").Value = '  -1.36%  '
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("D42").Value = '6.388'
$ws.Range("E42").Value = '  +6.56%  '
$ws.Range("D43").Value = '0.9067'
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '103.70'
$ws.Range("E45").Value = '  +0.33%  '
$ws.Range("D46").Value = '0.07971'
$ws.Range("E46").Value = '  +2.58%  '
$ws.Range("D47").Value = '1.990.41'
$ws.Range("E47").Value = '  +6.50%  '
$ws.Range("D48").Value = '65.28'
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("E49").Value = '  -2.99%  '
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '1.755'
$ws.Range("E51").Value = '  -0.81%  '
